# The source export re-sorted a block of observation rows (the underlying
# taxon sighting records for "Trossbygget, Dlr" got re-ordered), so on disk
# this shows up as the row contents being permuted across fixed row numbers:
#   rows 7 -> 8 -> 9 -> 7   (3-cycle)
#   rows 10 <-> 11          (swap)
#   rows 22 <-> 23           (swap)
#   rows 24 <-> 25           (swap)
# Row/record identity (Id, TaxonId, species names, coordinates, comments,
# etc.) moves with the record; only the row positions on the sheet change.
# Below we just write each destination row's new per-column values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows 7, 8, 9 (3-cycle: 7<-8, 8<-9, 9<-7) ---------------------------
$ws.Range("A7").Value = 131017563
$ws.Range("B7").Value = 79243
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("M7").Value = ""
$ws.Range("Q7").Value = 477226
$ws.Range("R7").Value = 6789084
$ws.Range("AC7").Value = "Rikligt i området"

$ws.Range("A8").Value = 131016886
$ws.Range("Q8").Value = 477116
$ws.Range("R8").Value = 6789167
$ws.Range("AC8").Value = ""

$ws.Range("A9").Value = 131017110
$ws.Range("B9").Value = 57884
$ws.Range("E9").Value = 100109
$ws.Range("F9").Value = "Tretåig hackspett"
$ws.Range("G9").Value = "Picoides tridactylus"
$ws.Range("H9").Value = "(Linnaeus, 1758)"
$ws.Range("M9").Value = "färska spår"
$ws.Range("Q9").Value = 477185
$ws.Range("R9").Value = 6789174

# --- rows 10, 11 (swap) --------------------------------------------------
$ws.Range("A10").Value = 131016935
$ws.Range("J10").Value = ""
$ws.Range("K10").Value = ""
$ws.Range("N10").Value = ""
$ws.Range("Q10").Value = 477129
$ws.Range("R10").Value = 6789191
$ws.Range("Z10").Value = "12:48"
$ws.Range("AB10").Value = "12:48"
$ws.Range("AF10").Value = ""

$ws.Range("A11").Value = 131023058
$ws.Range("J11").Value = ""
$ws.Range("K11").Value = ""
$ws.Range("N11").Value = ""
$ws.Range("Q11").Value = 477070
$ws.Range("R11").Value = 6788943
$ws.Range("Z11").Value = ""
$ws.Range("AB11").Value = ""
$ws.Range("AF11").Value = ""

# --- rows 22, 23 (swap) ---------------------------------------------------
$ws.Range("A22").Value = 131022789
$ws.Range("J22").Value = ""
$ws.Range("K22").Value = ""
$ws.Range("N22").Value = ""
$ws.Range("Q22").Value = 477076
$ws.Range("R22").Value = 6788997
$ws.Range("Z22").Value = ""
$ws.Range("AB22").Value = ""
$ws.Range("AF22").Value = ""

$ws.Range("A23").Value = 131016974
$ws.Range("J23").Value = ""
$ws.Range("K23").Value = ""
$ws.Range("N23").Value = ""
$ws.Range("Q23").Value = 477179
$ws.Range("R23").Value = 6789184
$ws.Range("Z23").Value = "12:48"
$ws.Range("AB23").Value = "12:48"
$ws.Range("AF23").Value = ""

# --- rows 24, 25 (swap) ---------------------------------------------------
$ws.Range("A24").Value = 131022847
$ws.Range("B24").Value = 57884
$ws.Range("E24").Value = 100109
$ws.Range("F24").Value = "Tretåig hackspett"
$ws.Range("G24").Value = "Picoides tridactylus"
$ws.Range("H24").Value = "(Linnaeus, 1758)"
$ws.Range("K24").Value = ""
$ws.Range("L24").Value = ""
$ws.Range("M24").Value = "färska spår"
$ws.Range("N24").Value = ""
$ws.Range("Q24").Value = 477122
$ws.Range("R24").Value = 6788910
$ws.Range("Z24").Value = ""
$ws.Range("AB24").Value = ""
$ws.Range("AC24").Value = ""

$ws.Range("A25").Value = 131017116
$ws.Range("B25").Value = 79243
$ws.Range("E25").Value = 6425
$ws.Range("F25").Value = "Garnlav"
$ws.Range("G25").Value = "Alectoria sarmentosa"
$ws.Range("H25").Value = "(Ach.) Ach."
$ws.Range("K25").Value = ""
$ws.Range("L25").Value = ""
$ws.Range("M25").Value = ""
$ws.Range("N25").Value = ""
$ws.Range("Q25").Value = 477185
$ws.Range("R25").Value = 6789174
$ws.Range("Z25").Value = "12:48"
$ws.Range("AB25").Value = "12:48"
$ws.Range("AC25").Value = "Rikligt i området"
